# Logged Week 16 and performed season sim from Week 17
#
# Updates the cumulative stat lines on both the "Rushing" and "Receiving"
# sheets, and adds a newly-appearing player to the bottom of each table
# (N.Foles does not get a new row on Rushing here - he's inserted into
# the existing table; C.Kmet / J.Horsted are the genuinely new rows).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Rushing"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Rushing")

# Rows 2 & 3 (A.Dalton, J.Fields) are unchanged.

# Existing rows 4-12 get updated names/values (table shifts down to make
# room for N.Foles, and a new player C.Kmet is appended as row 13).
$ws1.Range("B4").Value = "N.Foles"
$ws1.Range("C4").Value = 1
$ws1.Range("D4").Value = 0
$ws1.Range("E4").Value = 2
$ws1.Range("F4").Value = 2

$ws1.Range("B5").Value = "D.Montgomery"
$ws1.Range("C5").Value = 139
$ws1.Range("D5").Value = 71
$ws1.Range("E5").Value = 26
$ws1.Range("F5").Value = 29

$ws1.Range("B6").Value = "D.Williams"
$ws1.Range("C6").Value = 9
$ws1.Range("D6").Value = 10
$ws1.Range("E6").Value = 1
$ws1.Range("F6").Value = 0

$ws1.Range("B7").Value = "K.Herbert"
$ws1.Range("C7").Value = 16
$ws1.Range("D7").Value = 9
$ws1.Range("E7").Value = 2
$ws1.Range("F7").Value = 1

$ws1.Range("B8").Value = "R.Nall"
$ws1.Range("C8").Value = 1
$ws1.Range("D8").Value = 0
$ws1.Range("E8").Value = 0
$ws1.Range("F8").Value = 0

$ws1.Range("B9").Value = "A.Robinson"
$ws1.Range("C9").Value = 1
$ws1.Range("D9").Value = 0
$ws1.Range("E9").Value = 1
$ws1.Range("F9").Value = 0

$ws1.Range("B10").Value = "D.Mooney"
$ws1.Range("C10").Value = 4
$ws1.Range("D10").Value = 2
$ws1.Range("E10").Value = 0
$ws1.Range("F10").Value = 3

$ws1.Range("B11").Value = "M.Goodwin"
$ws1.Range("C11").Value = 2
$ws1.Range("D11").Value = 0
$ws1.Range("E11").Value = 0
$ws1.Range("F11").Value = 1

$ws1.Range("B12").Value = "J.Grant"
$ws1.Range("C12").Value = 3
$ws1.Range("D12").Value = 1
$ws1.Range("E12").Value = 1
$ws1.Range("F12").Value = 0

# New row 13 - C.Kmet. Copy formatting from the row above first so the
# styled (bold/border/centered) A column cell matches the rest of the
# table, then overwrite the values.
$ws1.Range("A12:F12").Copy($ws1.Range("A13:F13"))
$ws1.Range("A13").Value = 11
$ws1.Range("B13").Value = "C.Kmet"
$ws1.Range("C13").Value = 1
$ws1.Range("D13").Value = 0
$ws1.Range("E13").Value = 0
$ws1.Range("F13").Value = 1

# ---------------------------------------------------------------------
# Sheet "Receiving"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Receiving")

$ws2.Range("B2").Value = "D.Montgomery"
$ws2.Range("C2").Value = 55
$ws2.Range("D2").Value = 47
$ws2.Range("E2").Value = 1
$ws2.Range("F2").Value = 1
$ws2.Range("G2").Value = 9
$ws2.Range("H2").Value = 8

$ws2.Range("B3").Value = "D.Williams"
$ws2.Range("C3").Value = 11
$ws2.Range("D3").Value = 10
$ws2.Range("E3").Value = 2
$ws2.Range("F3").Value = 1
$ws2.Range("G3").Value = 2
$ws2.Range("H3").Value = 2

$ws2.Range("B4").Value = "K.Herbert"
$ws2.Range("C4").Value = 9
$ws2.Range("D4").Value = 8
$ws2.Range("E4").Value = 1
$ws2.Range("F4").Value = 1
$ws2.Range("G4").Value = 0
$ws2.Range("H4").Value = 0

$ws2.Range("B5").Value = "A.Robinson"
$ws2.Range("C5").Value = 41
$ws2.Range("D5").Value = 27
$ws2.Range("E5").Value = 16
$ws2.Range("F5").Value = 5
$ws2.Range("G5").Value = 6
$ws2.Range("H5").Value = 5

$ws2.Range("B6").Value = "D.Mooney"
$ws2.Range("C6").Value = 82
$ws2.Range("D6").Value = 48
$ws2.Range("E6").Value = 29
$ws2.Range("F6").Value = 14
$ws2.Range("G6").Value = 8
$ws2.Range("H6").Value = 4

$ws2.Range("B7").Value = "M.Goodwin"
$ws2.Range("C7").Value = 24
$ws2.Range("D7").Value = 13
$ws2.Range("E7").Value = 11
$ws2.Range("F7").Value = 5
$ws2.Range("G7").Value = 1
$ws2.Range("H7").Value = 1

$ws2.Range("B8").Value = "D.Byrd"
$ws2.Range("C8").Value = 25
$ws2.Range("D8").Value = 18
$ws2.Range("E8").Value = 5
$ws2.Range("F8").Value = 2
$ws2.Range("G8").Value = 4
$ws2.Range("H8").Value = 2

$ws2.Range("B9").Value = "J.Grant"
$ws2.Range("C9").Value = 11
$ws2.Range("D9").Value = 7
$ws2.Range("E9").Value = 4
$ws2.Range("F9").Value = 2
$ws2.Range("G9").Value = 1
$ws2.Range("H9").Value = 1

$ws2.Range("B10").Value = "D.Newsome"
$ws2.Range("C10").Value = 2
$ws2.Range("D10").Value = 1
$ws2.Range("E10").Value = 0
$ws2.Range("F10").Value = 0
$ws2.Range("G10").Value = 0
$ws2.Range("H10").Value = 0

$ws2.Range("B11").Value = "C.Kmet"
$ws2.Range("C11").Value = 73
$ws2.Range("D11").Value = 47
$ws2.Range("E11").Value = 10
$ws2.Range("F11").Value = 6
$ws2.Range("G11").Value = 12
$ws2.Range("H11").Value = 5

$ws2.Range("B12").Value = "J.Graham"
$ws2.Range("C12").Value = 18
$ws2.Range("D12").Value = 10
$ws2.Range("E12").Value = 2
$ws2.Range("F12").Value = 2
$ws2.Range("G12").Value = 8
$ws2.Range("H12").Value = 4

$ws2.Range("B13").Value = "J.James"
$ws2.Range("C13").Value = 7
$ws2.Range("D13").Value = 6
$ws2.Range("E13").Value = 0
$ws2.Range("F13").Value = 0
$ws2.Range("G13").Value = 2
$ws2.Range("H13").Value = 1

# New row 14 - J.Horsted. Copy formatting from the row above first.
$ws2.Range("A13:H13").Copy($ws2.Range("A14:H14"))
$ws2.Range("A14").Value = 12
$ws2.Range("B14").Value = "J.Horsted"
$ws2.Range("C14").Value = 1
$ws2.Range("D14").Value = 1
$ws2.Range("E14").Value = 0
$ws2.Range("F14").Value = 0
$ws2.Range("G14").Value = 1
$ws2.Range("H14").Value = 1
